$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AEE")

# Row 4: Inventory
$ws.Range("B4").Value = 521000000.0
$ws.Range("C4").Value = 557000000.0
$ws.Range("D4").Value = 514000000.0
$ws.Range("E4").Value = 471000000.0
$ws.Range("F4").Value = 494000000.0

# Row 14: Accounts Payable
$ws.Range("B14").Value = 958000000.0
$ws.Range("C14").Value = 640000000.0
$ws.Range("D14").Value = 616000000.0
$ws.Range("E14").Value = 544000000.0
$ws.Range("F14").Value = 874000000.0

# Row 22: Long Term Tax Liability (Deferred)
$ws.Range("B22").Value = 3211000000.0
$ws.Range("C22").Value = 3151000000.0
$ws.Range("D22").Value = 3034000000.0
$ws.Range("E22").Value = 2948000000.0
$ws.Range("F22").Value = 2919000000.0
